$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Reduce the "accuracy" (precision) of row 5's numeric data to 2 decimals ---
$row5 = @{
    "B5"  = 16.33
    "C5"  = 11.93
    "D5"  = 1.1
    "E5"  = 35.48
    "F5"  = 28.95
    "G5"  = 12.85
    "H5"  = 48.14
    "I5"  = 19.78
    "J5"  = 8.720000000000001
    "K5"  = 12.88
    "L5"  = 14.24
    "M5"  = 14.98
    "N5"  = 4.11
    "O5"  = 12.78
    "P5"  = 18.13
    "Q5"  = 10.86
    "R5"  = 0.8100000000000001
    "S5"  = 0.71
    "T5"  = 187.24
    "U5"  = 35.71
    "V5"  = 11.8
    "W5"  = 23.91
    "X5"  = 12.53
    "Y5"  = 1.99
    "Z5"  = 23.61
    "AA5" = 10.42
    "AB5" = 9.300000000000001
    "AC5" = 10.92
    "AD5" = 14.9
    "AE5" = 0.5600000000000001
    "AF5" = 43.61
    "AG5" = 6.6
    "AH5" = 14.75
}

foreach ($addr in $row5.Keys) {
    $ws.Range($addr).Value = $row5[$addr]
}

# --- 2. Delete row 6 (the extra sample row) entirely, shrinking the used range ---
$ws.Rows(6).Delete()

# --- 3. Narrow column AH (34) to match the other "7"-wide columns ---
$ws.Columns(34).ColumnWidth = $ws.Columns(4).ColumnWidth
